# Applies the changes described by the diff:
#  - Rename sheet "Işler" (dotless capital I) to "İşler" (dotted capital İ)
#  - Make "İşler" the active/selected tab (workbookView activeTab = index 2)
#  - On the "İşler" sheet, select cell H23 and mark the sheet as the tab-selected one
#  - On the "Birim" sheet, clear the tab-selected / previous selection state (A3 stays default)

$wb = $excel.ActiveWorkbook

# Rename the third sheet "Işler" -> "İşler"
$wsIsler = $wb.Worksheets.Item("Işler")
$wsIsler.Name = "İşler"

$wsBirim = $wb.Worksheets.Item("Birim")

# Update selection on the "İşler" sheet to H23 and activate it (this also
# makes it the active tab, matching activeTab="2" / tabSelected="1").
$wsIsler.Activate()
$wsIsler.Range("H23").Select()

# Touch the Birim sheet's selection back to its original default (A3) so it
# no longer carries the tab-selected marker.
$wsBirim.Range("A3").Select()

# Re-activate the İşler sheet last so it ends up as the saved active tab.
$wsIsler.Activate()
